# Regenerate orders with updated distance/size codes.
# The underlying experiment data encodes Distance (D51/D64/D80) and
# Size (S20/S25/S30) tokens inside several text columns (Condition,
# Filename_Left, Filename_Right, Distance, Size). This script renames
# those tokens workbook-wide:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# using literal, whole-text substring replacement (exactly like the
# Excel "Replace All" feature) so every occurrence - whether it is a
# standalone code like "D51" or embedded in a longer string like
# "Face08_D51_S25" or "Face08_D51_S25_l.png" - gets updated consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange

# Map of literal find -> replace strings. Order doesn't matter here
# because the four tokens (D51, D64, D80, S30) are mutually exclusive
# substrings, so applying them in any sequence yields the same result.
$replacements = [ordered]@{
    'D51' = 'D55'
    'D64' = 'D69'
    'D80' = 'D86'
    'S30' = 'S31'
}

foreach ($find in $replacements.Keys) {
    $replace = $replacements[$find]
    # xlWhole=1/xlPart=2 lookAt, we want part-match (2) since tokens are
    # embedded inside longer strings such as Face08_D51_S25.
    # Replace signature: Find, Replacement, LookAt, SearchOrder, MatchCase,
    # MatchByte, SearchFormat, ReplaceFormat
    $usedRange.Replace($find, $replace, 2, 1, $false, $false) | Out-Null
}
